# Auto-generated edit script applying the Adamantoise_Profits diff
# Updates market-price-derived columns (H,I,J,K,L,M,N) across several sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4086.7742
$ws.Range("J17").Value = 4086.7742
$ws.Range("L17").Value = 12260.3226
$ws.Range("N17").Value = -12596.3226
$ws.Range("H46").Value = 4990
$ws.Range("I46").Value = 4990
$ws.Range("K46").Value = 14970
$ws.Range("M46").Value = -14851
$ws.Range("H60").Value = 4990
$ws.Range("I60").Value = 4990
$ws.Range("K60").Value = 14970
$ws.Range("M60").Value = -14486

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2244.9333
$ws.Range("I2").Value = 1350
$ws.Range("J2").Value = 2382.6155
$ws.Range("K2").Value = 1350
$ws.Range("L2").Value = 2382.6155
$ws.Range("M2").Value = -1237
$ws.Range("N2").Value = -2608.6155
$ws.Range("H45").Value = 4629.1665
$ws.Range("I45").Value = 5161.4546
$ws.Range("K45").Value = 5161.4546
$ws.Range("M45").Value = -4784.4546
$ws.Range("H63").Value = 3661
$ws.Range("I63").Value = 2189.8
$ws.Range("J63").Value = 5500
$ws.Range("K63").Value = 2189.8
$ws.Range("L63").Value = 5500
$ws.Range("M63").Value = -1503.8
$ws.Range("N63").Value = -6872
$ws.Range("H66").Value = 3661
$ws.Range("I66").Value = 2189.8
$ws.Range("J66").Value = 5500
$ws.Range("K66").Value = 10949
$ws.Range("L66").Value = 27500
$ws.Range("M66").Value = -7517
$ws.Range("N66").Value = -34364
$ws.Range("H102").Value = 1843.5625
$ws.Range("I102").Value = 1499.8182
$ws.Range("K102").Value = 1499.8182
$ws.Range("M102").Value = 122.1818000000001
$ws.Range("H116").Value = 2244.9333
$ws.Range("I116").Value = 1350
$ws.Range("J116").Value = 2382.6155
$ws.Range("K116").Value = 1350
$ws.Range("L116").Value = 2382.6155
$ws.Range("M116").Value = 944
$ws.Range("N116").Value = -6970.6155
$ws.Range("H122").Value = 4364.7
$ws.Range("I122").Value = 3749.5789
$ws.Range("K122").Value = 11248.7367
$ws.Range("M122").Value = -8798.736699999999
$ws.Range("H132").Value = 2558.5676
$ws.Range("I132").Value = 2558.2942
$ws.Range("K132").Value = 7674.882599999999
$ws.Range("M132").Value = -5144.882599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2244.9333
$ws.Range("I3").Value = 1350
$ws.Range("J3").Value = 2382.6155
$ws.Range("K3").Value = 1350
$ws.Range("L3").Value = 2382.6155
$ws.Range("M3").Value = -1236
$ws.Range("N3").Value = -2610.6155
$ws.Range("H11").Value = 1938.3334
$ws.Range("I11").Value = 10
$ws.Range("K11").Value = 10
$ws.Range("M11").Value = 130
$ws.Range("H20").Value = 6998.5
$ws.Range("I20").Value = 6998.5
$ws.Range("K20").Value = 6998.5
$ws.Range("M20").Value = -6751.5
$ws.Range("H86").Value = 3309.5454
$ws.Range("I86").Value = 2963.2917
$ws.Range("K86").Value = 2963.2917
$ws.Range("M86").Value = -1840.2917
$ws.Range("H89").Value = 3309.5454
$ws.Range("I89").Value = 2963.2917
$ws.Range("K89").Value = 14816.4585
$ws.Range("M89").Value = -9200.458500000001
$ws.Range("H99").Value = 1731.6875
$ws.Range("I99").Value = 1208.2307
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 1208.2307
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = 289.7692999999999
$ws.Range("N99").Value = -6996
$ws.Range("H105").Value = 2721.375
$ws.Range("I105").Value = 2378.75
$ws.Range("K105").Value = 2378.75
$ws.Range("M105").Value = -631.75
$ws.Range("H134").Value = 2901309.2
$ws.Range("I134").Value = 3177148.2
$ws.Range("K134").Value = 9531444.600000001
$ws.Range("M134").Value = -9528909.600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8716.286
$ws.Range("I31").Value = 3995
$ws.Range("J31").Value = 8952.35
$ws.Range("K31").Value = 3995
$ws.Range("L31").Value = 8952.35
$ws.Range("M31").Value = -3700
$ws.Range("N31").Value = -9542.35
$ws.Range("H34").Value = 8716.286
$ws.Range("I34").Value = 3995
$ws.Range("J34").Value = 8952.35
$ws.Range("K34").Value = 3995
$ws.Range("L34").Value = 8952.35
$ws.Range("M34").Value = -3793
$ws.Range("N34").Value = -9356.35
$ws.Range("H68").Value = 49962.145
$ws.Range("J68").Value = 49962.145
$ws.Range("L68").Value = 49962.145
$ws.Range("N68").Value = -51460.145
$ws.Range("H71").Value = 49962.145
$ws.Range("J71").Value = 49962.145
$ws.Range("L71").Value = 149886.435
$ws.Range("N71").Value = -157374.435
$ws.Range("H99").Value = 3093.7144
$ws.Range("I99").Value = 2971.2
$ws.Range("K99").Value = 2971.2
$ws.Range("M99").Value = -1473.2
$ws.Range("H105").Value = 3104.8572
$ws.Range("I105").Value = 2951.6365
$ws.Range("K105").Value = 2951.6365
$ws.Range("M105").Value = -1204.6365
$ws.Range("H107").Value = 32235.781
$ws.Range("I107").Value = 55903.555
$ws.Range("K107").Value = 55903.555
$ws.Range("M107").Value = -53983.555
$ws.Range("H122").Value = 5239.4
$ws.Range("I122").Value = 4346.4443
$ws.Range("K122").Value = 13039.3329
$ws.Range("M122").Value = -10589.3329
$ws.Range("H126").Value = 3093.7144
$ws.Range("I126").Value = 2971.2
$ws.Range("K126").Value = 8913.599999999999
$ws.Range("M126").Value = -6443.599999999999
$ws.Range("H132").Value = 2612.9524
$ws.Range("J132").Value = 4675.6665
$ws.Range("L132").Value = 14026.9995
$ws.Range("N132").Value = -19086.9995
$ws.Range("H134").Value = 3136.8
$ws.Range("I134").Value = 3153.111
$ws.Range("J134").Value = 2990
$ws.Range("K134").Value = 9459.332999999999
$ws.Range("L134").Value = 8970
$ws.Range("M134").Value = -6924.332999999999
$ws.Range("N134").Value = -14040

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2397.75
$ws.Range("I80").Value = 2071.875
$ws.Range("J80").Value = 3049.5
$ws.Range("K80").Value = 2071.875
$ws.Range("L80").Value = 3049.5
$ws.Range("M80").Value = -1073.875
$ws.Range("N80").Value = -5045.5
$ws.Range("H83").Value = 2397.75
$ws.Range("I83").Value = 2071.875
$ws.Range("J83").Value = 3049.5
$ws.Range("K83").Value = 10359.375
$ws.Range("L83").Value = 15247.5
$ws.Range("M83").Value = -5367.375
$ws.Range("N83").Value = -25231.5
$ws.Range("H97").Value = 531.1667
$ws.Range("I97").Value = 414.76923
$ws.Range("J97").Value = 833.8
$ws.Range("K97").Value = 414.76923
$ws.Range("L97").Value = 833.8
$ws.Range("M97").Value = 81.23077000000001
$ws.Range("N97").Value = -1825.8
$ws.Range("H113").Value = 25867.297
$ws.Range("I113").Value = 5496.25
$ws.Range("J113").Value = 55497.91
$ws.Range("K113").Value = 5496.25
$ws.Range("L113").Value = 55497.91
$ws.Range("M113").Value = -3326.25
$ws.Range("N113").Value = -59837.91
$ws.Range("H122").Value = 1777.1333
$ws.Range("I122").Value = 1653.6364
$ws.Range("J122").Value = 2116.75
$ws.Range("K122").Value = 4960.9092
$ws.Range("L122").Value = 6350.25
$ws.Range("M122").Value = -2510.9092
$ws.Range("N122").Value = -11250.25
$ws.Range("H126").Value = 2757
$ws.Range("J126").Value = 3009.3333
$ws.Range("L126").Value = 9027.999899999999
$ws.Range("N126").Value = -13967.9999
$ws.Range("H132").Value = 1970.8
$ws.Range("I132").Value = 2083.5881
$ws.Range("K132").Value = 6250.7643
$ws.Range("M132").Value = -3720.7643

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("H55").Value = 847.02563
$ws.Range("I55").Value = 791.93335
$ws.Range("K55").Value = 791.93335
$ws.Range("M55").Value = -618.93335
$ws.Range("H82").Value = 2201
$ws.Range("I82").Value = 1690.5714
$ws.Range("K82").Value = 1690.5714
$ws.Range("M82").Value = -1329.5714
$ws.Range("H85").Value = 2201
$ws.Range("I85").Value = 1690.5714
$ws.Range("K85").Value = 1690.5714
$ws.Range("M85").Value = -442.5714
$ws.Range("H93").Value = 27779214
$ws.Range("I93").Value = 55556764
$ws.Range("K93").Value = 55556764
$ws.Range("M93").Value = -55555516
$ws.Range("H100").Value = 1332.3334
$ws.Range("I100").Value = 998.5
$ws.Range("K100").Value = 998.5
$ws.Range("M100").Value = -457.5
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 10634.786
$ws.Range("I136").Value = 14957.333
$ws.Range("J136").Value = 9455.909
$ws.Range("K136").Value = 44871.999
$ws.Range("L136").Value = 28367.727
$ws.Range("M136").Value = -42321.999
$ws.Range("N136").Value = -33467.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 493.36365
$ws.Range("I100").Value = 131.14285
$ws.Range("J100").Value = 1127.25
$ws.Range("K100").Value = 262.2857
$ws.Range("L100").Value = 2254.5
$ws.Range("M100").Value = 278.7143
$ws.Range("N100").Value = -3336.5
$ws.Range("H107").Value = 442.88235
$ws.Range("I107").Value = 397.4
$ws.Range("K107").Value = 1192.2
$ws.Range("M107").Value = 727.8000000000002
$ws.Range("H113").Value = 410.5
$ws.Range("J113").Value = 601
$ws.Range("L113").Value = 1803
$ws.Range("N113").Value = -6143
$ws.Range("H122").Value = 5338.033
$ws.Range("I122").Value = 5004.857
$ws.Range("K122").Value = 15014.571
$ws.Range("M122").Value = -12564.571
$ws.Range("H132").Value = 3050.1853
$ws.Range("J132").Value = 3579.4
$ws.Range("L132").Value = 10738.2
$ws.Range("N132").Value = -15798.2
